$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2215.2856
$ws.Range("I18").Value = 2621
$ws.Range("K18").Value = 2621
$ws.Range("M18").Value = -2337
$ws.Range("H28").Value = 535
$ws.Range("I28").Value = 513.86365
$ws.Range("K28").Value = 513.86365
$ws.Range("M28").Value = -28.86365000000001
$ws.Range("H43").Value = 3474.75
$ws.Range("J43").Value = 3474.75
$ws.Range("L43").Value = 3474.75
$ws.Range("N43").Value = -3612.75
$ws.Range("H104").Value = 451
$ws.Range("I104").Value = 451
$ws.Range("K104").Value = 1353
$ws.Range("M104").Value = 394
$ws.Range("H112").Value = 3147.75
$ws.Range("J112").Value = 4197
$ws.Range("L112").Value = 12591
$ws.Range("N112").Value = -14807
$ws.Range("H116").Value = 25132.406
$ws.Range("J116").Value = 28224.85
$ws.Range("L116").Value = 28224.85
$ws.Range("N116").Value = -35108.85
$ws.Range("H132").Value = 9030.243
$ws.Range("I132").Value = 9267.194
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 27801.582
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -25271.582
$ws.Range("N132").Value = -6560
$ws.Range("H137").Value = 7834.846
$ws.Range("I137").Value = 2456
$ws.Range("K137").Value = 7368
$ws.Range("M137").Value = -4818
$ws.Range("H138").Value = 6565.75
$ws.Range("I138").Value = 10245.272
$ws.Range("J138").Value = 2068.5557
$ws.Range("K138").Value = 30735.816
$ws.Range("L138").Value = 6205.6671
$ws.Range("M138").Value = -25595.816
$ws.Range("N138").Value = -16485.6671
$ws.Range("H141").Value = 1344.0264
$ws.Range("I141").Value = 1176.3715
$ws.Range("J141").Value = 3300
$ws.Range("K141").Value = 3529.1145
$ws.Range("L141").Value = 9900
$ws.Range("M141").Value = 1650.8855
$ws.Range("N141").Value = -20260

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 3690
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 5683.3335
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 5683.3335
$ws.Range("M22").Value = -401
$ws.Range("N22").Value = -6281.3335
$ws.Range("H88").Value = 1326.2609
$ws.Range("I88").Value = 698.8333
$ws.Range("J88").Value = 1547.7059
$ws.Range("K88").Value = 698.8333
$ws.Range("L88").Value = 1547.7059
$ws.Range("M88").Value = -292.8333
$ws.Range("N88").Value = -2359.7059
$ws.Range("H91").Value = 1326.2609
$ws.Range("I91").Value = 698.8333
$ws.Range("J91").Value = 1547.7059
$ws.Range("K91").Value = 698.8333
$ws.Range("L91").Value = 1547.7059
$ws.Range("M91").Value = 705.1667
$ws.Range("N91").Value = -4355.7059
$ws.Range("H102").Value = 4365.696
$ws.Range("I102").Value = 3356.2222
$ws.Range("K102").Value = 3356.2222
$ws.Range("M102").Value = -1734.2222

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 31110.756
$ws.Range("I20").Value = 38622.965
$ws.Range("K20").Value = 38622.965
$ws.Range("M20").Value = -38375.965
$ws.Range("H22").Value = 749.75
$ws.Range("I22").Value = 749.75
$ws.Range("K22").Value = 749.75
$ws.Range("M22").Value = -576.75
$ws.Range("H62").Value = 50001
$ws.Range("J62").Value = 50001
$ws.Range("L62").Value = 50001
$ws.Range("N62").Value = -51373
$ws.Range("H65").Value = 50001
$ws.Range("J65").Value = 50001
$ws.Range("L65").Value = 150003
$ws.Range("N65").Value = -156867
$ws.Range("H94").Value = 3107
$ws.Range("I94").Value = 2408.8462
$ws.Range("J94").Value = 7645
$ws.Range("K94").Value = 2408.8462
$ws.Range("L94").Value = 7645
$ws.Range("M94").Value = -1957.8462
$ws.Range("N94").Value = -8547
$ws.Range("H97").Value = 31562.5
$ws.Range("I97").Value = 22356.5
$ws.Range("K97").Value = 22356.5
$ws.Range("M97").Value = -21365.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4243.636
$ws.Range("I31").Value = 4306.727
$ws.Range("K31").Value = 4306.727
$ws.Range("M31").Value = -4011.727
$ws.Range("H34").Value = 4243.636
$ws.Range("I34").Value = 4306.727
$ws.Range("K34").Value = 4306.727
$ws.Range("M34").Value = -4104.727
$ws.Range("H58").Value = 17616
$ws.Range("I58").Value = 9408
$ws.Range("J58").Value = 22540.8
$ws.Range("K58").Value = 9408
$ws.Range("L58").Value = 22540.8
$ws.Range("M58").Value = -9205
$ws.Range("N58").Value = -22946.8
$ws.Range("H99").Value = 44179.6
$ws.Range("J99").Value = 3500
$ws.Range("L99").Value = 3500
$ws.Range("N99").Value = -6496
$ws.Range("H126").Value = 44179.6
$ws.Range("J126").Value = 3500
$ws.Range("L126").Value = 10500
$ws.Range("N126").Value = -15440
$ws.Range("H136").Value = 17616
$ws.Range("I136").Value = 9408
$ws.Range("J136").Value = 22540.8
$ws.Range("K136").Value = 28224
$ws.Range("L136").Value = 67622.39999999999
$ws.Range("M136").Value = -25674
$ws.Range("N136").Value = -72722.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 3769.1428
$ws.Range("I32").Value = 3463.3333
$ws.Range("J32").Value = 3998.5
$ws.Range("K32").Value = 10389.9999
$ws.Range("L32").Value = 11995.5
$ws.Range("M32").Value = -10106.9999
$ws.Range("N32").Value = -12561.5
$ws.Range("H97").Value = 291.75
$ws.Range("J97").Value = 297.33334
$ws.Range("L97").Value = 892.0000200000001
$ws.Range("N97").Value = -1884.00002
$ws.Range("H137").Value = 2241
$ws.Range("I137").Value = 2045.5555
$ws.Range("K137").Value = 6136.666499999999
$ws.Range("M137").Value = -1036.666499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2577.96
$ws.Range("I102").Value = 2601.4583
$ws.Range("K102").Value = 2601.4583
$ws.Range("M102").Value = -979.4582999999998
$ws.Range("H113").Value = 1647.3667
$ws.Range("I113").Value = 1544.4348
$ws.Range("J113").Value = 1985.5714
$ws.Range("K113").Value = 1544.4348
$ws.Range("L113").Value = 1985.5714
$ws.Range("M113").Value = 625.5652
$ws.Range("N113").Value = -6325.5714
$ws.Range("H126").Value = 7859.4546
$ws.Range("I126").Value = 7471.625
$ws.Range("J126").Value = 8893.666999999999
$ws.Range("K126").Value = 22414.875
$ws.Range("L126").Value = 26681.001
$ws.Range("M126").Value = -19944.875
$ws.Range("N126").Value = -31621.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1452.25
$ws.Range("I9").Value = 231.14285
$ws.Range("K9").Value = 231.14285
$ws.Range("M9").Value = -7.14285000000001
$ws.Range("H93").Value = 1945
$ws.Range("I93").Value = 1320.762
$ws.Range("K93").Value = 1320.762
$ws.Range("M93").Value = -72.76199999999994
$ws.Range("H136").Value = 7001
$ws.Range("I136").Value = 6999
$ws.Range("J136").Value = 7003
$ws.Range("K136").Value = 20997
$ws.Range("L136").Value = 21009
$ws.Range("M136").Value = -18447
$ws.Range("N136").Value = -26109

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H81").Value = 3794.2
$ws.Range("I81").Value = 3097.353
$ws.Range("J81").Value = 5275
$ws.Range("K81").Value = 6194.706
$ws.Range("L81").Value = 10550
$ws.Range("M81").Value = -5133.706
$ws.Range("N81").Value = -12672
$ws.Range("H84").Value = 3794.2
$ws.Range("I84").Value = 3097.353
$ws.Range("J84").Value = 5275
$ws.Range("K84").Value = 30973.53
$ws.Range("L84").Value = 52750
$ws.Range("M84").Value = -25669.53
$ws.Range("N84").Value = -63358
$ws.Range("H113").Value = 632.2963
$ws.Range("I113").Value = 541
$ws.Range("J113").Value = 1157.25
$ws.Range("K113").Value = 1623
$ws.Range("L113").Value = 3471.75
$ws.Range("M113").Value = 547
$ws.Range("N113").Value = -7811.75
$ws.Range("H122").Value = 42852.07
$ws.Range("I122").Value = 4320.706
$ws.Range("J122").Value = 97438.164
$ws.Range("K122").Value = 12962.118
$ws.Range("L122").Value = 292314.492
$ws.Range("M122").Value = -10512.118
$ws.Range("N122").Value = -297214.492
$ws.Range("H132").Value = 3474.8462
$ws.Range("I132").Value = 3308.3928
$ws.Range("J132").Value = 3898.5454
$ws.Range("K132").Value = 9925.178400000001
$ws.Range("L132").Value = 11695.6362
$ws.Range("M132").Value = -7395.178400000001
$ws.Range("N132").Value = -16755.6362
